$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (RM 8): F3 value 17.64 -> missing ---
$ws.Range("F3").ClearContents()

# --- Remove the "RM 232" row (original row 26) and the "SC 92" row (original row 28) ---
# Deleting row 26 first shifts "SC 92" (originally row 28) up to row 27.
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# After both deletions the sheet has rows 2-33 of data (32 records), dimension A1:F33.
# Remaining row edits (by the now-current row numbers):
#   row 26 = "SC 5"   -> D26 filled in as -13.8
#   row 27 = "SC 101" -> D27 becomes missing
#   row 33 = "SC 232" -> D33 filled in as -14.1, F33 filled in as 17.53

$ws.Range("D26").Value = -13.8
$ws.Range("D27").ClearContents()
$ws.Range("D33").Value = -14.1
$ws.Range("F33").Value = 17.53
